$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '28.514.94'
$ws.Cells.Item(2, 5).Value = '  +0.28%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.564.38'
$ws.Cells.Item(3, 5).Value = '  -2.03%  '

$ws.Cells.Item(4, 5).Value = '  +0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '211.63'
$ws.Cells.Item(5, 5).Value = '  -1.48%  '

$ws.Cells.Item(6, 5).Value = '  -1.10%  '

$ws.Cells.Item(7, 5).Value = '  +0.11%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '46.37'
$ws.Cells.Item(8, 5).Value = '  +5.34%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '24.08'
$ws.Cells.Item(9, 5).Value = '  +0.18%  '

$ws.Cells.Item(10, 5).Value = '  -1.76%  '

$ws.Cells.Item(11, 5).Value = '  -1.53%  '

$ws.Cells.Item(12, 5).Value = '  -0.26%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '1.787.29'
$ws.Cells.Item(13, 5).Value = '  -1.99%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '1.573.77'
$ws.Cells.Item(14, 5).Value = '  -1.64%  '

$ws.Cells.Item(15, 5).Value = '  -2.57%  '

$ws.Cells.Item(16, 2).Value = 'Polkadot'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '3.68'
$ws.Cells.Item(16, 5).Value = '  -3.27%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '28.509.49'
$ws.Cells.Item(17, 5).Value = '  +0.15%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '62.18'
$ws.Cells.Item(18, 5).Value = '  -1.79%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '228.96'
$ws.Cells.Item(19, 5).Value = '  -1.76%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.0₃0694'
$ws.Cells.Item(20, 5).Value = '  -2.46%  '

$ws.Cells.Item(21, 5).Value = '  -2.55%  '

$ws.Cells.Item(22, 5).Value = '  +0.10%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '3.88'
$ws.Cells.Item(23, 5).Value = '  -6.03%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '9.15'
$ws.Cells.Item(24, 5).Value = '  -3.03%  '

$ws.Cells.Item(25, 5).Value = '  +6.69%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '150.46'
$ws.Cells.Item(26, 5).Value = '  -1.46%  '

$ws.Cells.Item(27, 5).Value = '  -2.17%  '

$ws.Cells.Item(28, 5).Value = '  -2.79%  '

$ws.Cells.Item(29, 5).Value = '  -3.94%  '

$ws.Cells.Item(30, 5).Value = '  +0.08%  '

$ws.Cells.Item(31, 5).Value = '  -2.06%  '

$ws.Cells.Item(32, 5).Value = '  -4.07%  '

$ws.Cells.Item(33, 5).Value = '  -1.29%  '

$ws.Cells.Item(34, 5).Value = '  -2.62%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.391.92'
$ws.Cells.Item(35, 5).Value = '  -2.02%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.04'
$ws.Cells.Item(36, 5).Value = '  -1.37%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.55'
$ws.Cells.Item(37, 5).Value = '  -3.26%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.35'
$ws.Cells.Item(38, 5).Value = '  +0.83%  '

$ws.Cells.Item(40, 5).Value = '  -1.14%  '

$ws.Cells.Item(41, 5).Value = '  -1.98%  '

$ws.Cells.Item(42, 5).Value = '  +0.13%  '

$ws.Cells.Item(43, 5).Value = '  +2.94%  '

$ws.Cells.Item(44, 5).Value = '  -4.09%  '

$ws.Cells.Item(45, 5).Value = '  -4.38%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.976'
$ws.Cells.Item(46, 5).Value = '  -0.49%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '62.72'
$ws.Cells.Item(47, 5).Value = '  -3.32%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '1.700.48'
$ws.Cells.Item(48, 5).Value = '  -1.95%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '86.05'
$ws.Cells.Item(49, 5).Value = '  -1.81%  '

$ws.Cells.Item(50, 5).Value = '  -4.11%  '

$ws.Cells.Item(51, 5).Value = '  -0.19%  '
